# SCE_SV_SE.xlsx — "updated estimates for revised paper"
#
# Changes vs. the original table:
#   - A new (blank) column D is introduced to the right of the existing data.
#   - The row-1 header cells (B1/C1), which held the LaTeX column labels
#     ($\hat\lambda$ / $\gamma$), are reset to plain 0 placeholders, and the
#     new D1 cell gets the same 0 placeholder.
#   - The point estimates in column B (FE / FE+Disg / FE+Disg+Var rows) move
#     from 0.35 to the revised 0.36; column C (0.2) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new trailing column (D) used by the revised table layout.
$ws.Columns.Item(4).Insert()

# Row 1: header placeholders become numeric 0 (was text labels in B1/C1;
# D1 is the newly added column's placeholder).
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = 0

# Revised point estimates: 0.35 -> 0.36 for every data row in column B.
$ws.Range("B2").Value = 0.36
$ws.Range("B3").Value = 0.36
$ws.Range("B4").Value = 0.36
